$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("U2").Value = 3.8
$ws.Range("V2").Value = 1.27

# Row 4 updates
$ws.Range("G4").Value = 1.87
$ws.Range("M4").Value = 1.05
$ws.Range("O4").Value = 1.37
$ws.Range("X4").Value = 1.19

# Row 5 updates
$ws.Range("G5").Value = 1.47
$ws.Range("H5").Value = 3.7
$ws.Range("L5").Value = 6.5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 8.5
$ws.Range("O5").Value = 1.3
$ws.Range("X5").Value = 1.22
$ws.Range("AC5").Value = 5.5
$ws.Range("AG5").Value = 15
$ws.Range("AQ5").Value = 67

$wb.Save()
